# formatLaporanCuSemua.xlsx edit:
# - rename Sheet2 -> Sheet1
# - insert "periode" as new column B (shifting old B.. right)
# - insert "aset_likuid_tidak_menghasilkan" as new column L (shifting old L.. right)
# - append "tanggal buat" as new last column (AN)
# - add a few blank date-formatted cells in AM2:AM4 / AN2:AN4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename sheet
$ws.Name = "Sheet1"

# 2. Insert new column B for "periode" (old col B "lelaki biasa" etc. shift right to C..)
$ws.Columns(2).Insert()
$ws.Range("B1").Value = "periode"

# 3. Insert new column L for "aset_likuid_tidak_menghasilkan"
#    (after the insert above, "aset tidak menghasilkan" sits at column K,
#     "aktiva lancar" at column L - push it right to make room)
$ws.Columns(12).Insert()
$ws.Range("L1").Value = "aset_likuid_tidak_menghasilkan"

# 4. After the two inserts above, the old "periode" (originally AL) now
#    sits at AN1 and the old "tgl buat" (originally AM) now sits at AO1.
#    "periode" was already re-created fresh at B1, so AN1 is safe to
#    overwrite with the new "tanggal buat" header; the stale trailing
#    "tgl buat" column (now AO) is deleted so the sheet still ends at AN.
$ws.Range("AN1").Value = "tanggal buat"
$ws.Columns(41).Delete()

# 5. Add the new date-formatted (but empty) cells in rows 2-4.
#    Format one cell, then fan the format out with copy/paste-special so
#    every cell shares a single cellXfs entry (matches authored file,
#    which reuses style index 1 for all six cells) instead of each
#    Range.NumberFormat assignment minting its own style record.
$ws.Range("AM2").NumberFormat = "mm-dd-yy"
$ws.Range("AM2").Copy()
$ws.Range("AM3:AM4").PasteSpecial(-4122)
$ws.Range("AN2:AN4").PasteSpecial(-4122)

# 6. Match the author's final cursor position (B11) recorded in the sheet view.
[void]$ws.Range("B11").Select()
